$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New loading_percent values per column (rows 2-25), derived from the
# "case with 380 kV done" update.
$BValues = @(18.71113896771032, 18.17061707178904, 17.83442004836189, 17.69655126238573, 17.67361262419753, 17.83256390805282, 18.52578979285231, 19.84204279218102, 20.77186227699329, 21.18477064987499, 21.33953324844894, 21.3062755853303, 21.1975358645325, 21.13071752832144, 20.74466204673618, 20.50513962644961, 20.36643265327006, 20.3193121767173, 20.53073560778505, 21.22951980637182, 21.67684928086977, 21.43900350868357, 20.51916677524367, 19.49176824313264)
$CValues = @(9.711612252535442, 9.396928590263286, 9.196653604177982, 9.113341173616302, 9.099406981567117, 9.195536795587095, 9.604617440648138, 10.34786725434179, 10.85454388358351, 11.07589536486481, 11.15835941308749, 11.14066028138101, 11.0827071535817, 11.04703129362987, 10.8398903679636, 10.71044332914866, 10.63513087895515, 10.60948542251376, 10.72431230844315, 11.09976653438643, 11.33722202717566, 11.21122554492304, 10.7180449168666, 10.15349256044961)
$EValues = @(11.20705533867923, 11.24473097846779, 11.26957486677224, 11.28012917057682, 11.28190768892811, 11.2697154638653, 11.21969080287876, 11.13516778260468, 11.08134456805858, 11.0586562506748, 11.05032313608173, 11.0521063234654, 11.0579654993753, 11.06158807818005, 11.08286346066742, 11.09637537443476, 11.10431611925603, 11.10703374814946, 11.09491951095051, 11.05623750142718, 11.03246315652572, 11.04501405830609, 11.09557716932982, 11.15658021319915)
$FValues = @(16.86991607391233, 15.89584955866808, 15.26997757108491, 15.00819731993403, 14.96433081551593, 15.26647399323137, 16.53996406344768, 19.00274580682531, 20.67494806633232, 21.3917225636224, 21.65686569030329, 21.60004134736742, 21.4136618050453, 21.29868154950795, 20.62722412089977, 20.20408069617459, 19.95656407809808, 19.87204792380562, 20.24955283636157, 21.46857628470577, 22.22866616901552, 21.82633154475857, 20.22900810905287, 18.34778573295691)
$GValues = @(38.44006884876254, 38.42605159419219, 38.43721151655754, 38.44669488880991, 38.44856654438062, 38.43731947976798, 38.43111426944837, 38.57700029823251, 38.78168925773615, 38.89606688083302, 38.9424314140771, 38.93231040392367, 38.89982020695339, 38.8803162162125, 38.77464199889152, 38.71525898462899, 38.68310570032709, 38.67256307887537, 38.72137319746327, 38.90928061130271, 39.04987159819166, 38.97321223479829, 38.71860277285643, 38.52046493493215)
$HValues = @(16.47425620400155, 16.54236111961629, 16.58892091059538, 16.60908014058028, 16.61249897573836, 16.58918799337007, 16.49675019808315, 16.35342084168226, 16.27165902980875, 16.23966119317112, 16.22829845391358, 16.23071197318855, 16.23871122028322, 16.24370940599913, 16.27385532033324, 16.29368469883324, 16.30557870117321, 16.30968952009555, 16.29152320151249, 16.23634112553782, 16.20467529785113, 16.22117119636362, 16.29249887671686, 16.38809030694387)
$LValues = @(9.937090530883424, 9.914095803461628, 9.901887477985008, 9.897395957941125, 9.896679426066289, 9.901824942493375, 9.92876737056041, 9.996609423952947, 10.05537546942686, 10.08399061129173, 10.09509166201046, 10.09268915292678, 10.08489862105882, 10.08016105098692, 10.05354281176854, 10.0376916308367, 10.02875196496529, 10.02575579863207, 10.03936068303599, 10.08717973981081, 10.11997474737619, 10.10233221792709, 10.03860556404498, 9.976671444695553)

$startRow = 2
$columns = @(
    @{ Col = 2; Values = $BValues },
    @{ Col = 3; Values = $CValues },
    @{ Col = 5; Values = $EValues },
    @{ Col = 6; Values = $FValues },
    @{ Col = 7; Values = $GValues },
    @{ Col = 8; Values = $HValues },
    @{ Col = 12; Values = $LValues }
)

foreach ($entry in $columns) {
    $c = $entry.Col
    $vals = $entry.Values
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($startRow + $i, $c).Value = $vals[$i]
    }
}
